$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap the contents of columns F..V between two rows (everything
# except the leading A..E "match id / country / league / season / date"
# columns, which stay put).
# ---------------------------------------------------------------------------
function Swap-MatchData($row1, $row2) {
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($c in $cols) {
        $cell1 = $c + $row1
        $cell2 = $c + $row2
        $v1 = $ws.Range($cell1).Value()
        $v2 = $ws.Range($cell2).Value()
        $ws.Range($cell1).Value = $v2
        $ws.Range($cell2).Value = $v1
    }
}

# Re-ordered match rows (odds data refreshed / rows resorted by the scraper)
Swap-MatchData 106 107
Swap-MatchData 123 124
Swap-MatchData 126 127
Swap-MatchData 155 156

# ---------------------------------------------------------------------------
# Append the new match row (171) at the end of the sheet, copying the
# formatting of the last existing data row (170) first.
# ---------------------------------------------------------------------------
$ws.Range("A170:V170").Copy()
$ws.Range("A171:V171").PasteSpecial(-4122)

$ws.Range("A171").Value = 170
$ws.Range("B171").Value = "spain"
$ws.Range("C171").Value = "primera-rfef-group-1"
$ws.Range("D171").Value = "2023-2024"
$ws.Range("E171").Value = 45293.79166666666
$ws.Range("F171").Value = "Lugo"
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = "Celta Vigo B"
$ws.Range("I171").Value = 1
$ws.Range("J171").Value = 2.42
$ws.Range("K171").Value = "31/12/2024 19:12"
$ws.Range("L171").Value = 2.49
$ws.Range("M171").Value = "02/01/2024 18:29"
$ws.Range("N171").Value = 2.89
$ws.Range("O171").Value = "31/12/2024 19:12"
$ws.Range("P171").Value = 3.05
$ws.Range("Q171").Value = "02/01/2024 18:27"
$ws.Range("R171").Value = 2.98
$ws.Range("S171").Value = "31/12/2024 19:12"
$ws.Range("T171").Value = 3.05
$ws.Range("U171").Value = "02/01/2024 18:29"
$ws.Range("V171").Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-1/lugo-celta-vigo/8IsCzT2D/"
